# Add a new "E" column that duplicates the numeric value shown in column D
# (the "header_col_0" / year column) into a new "year" column for the
# worksheets that still only have that value in D.
#
# Sheets 1, 10, 14, 15 and 16 (1-based tab order) each have a header in
# E1 ("year") but are still missing the per-row data in column E; every
# other data row (2..last) needs E<row> = numeric value of D<row>.

$wb = $excel.ActiveWorkbook

$targetSheets = @(1, 10, 14, 15, 16)

foreach ($sheetIndex in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # Determine the last used row from the sheet's UsedRange (some rows have
    # gaps in column A, so we can't rely on a single column for End(xlUp)).
    $usedRange = $ws.UsedRange
    $firstRow = $usedRange.Row
    $lastRow = $firstRow + $usedRange.Rows.Count - 1

    for ($r = 2; $r -le $lastRow; $r++) {
        $dCell = $ws.Cells.Item($r, 4)
        $eCell = $ws.Cells.Item($r, 5)

        # Value2 gives the underlying value (numeric for year cells, whether
        # the source cell is stored as a shared string or a plain number).
        $yearValue = $dCell.Value2

        if ($null -ne $yearValue -and [string]$yearValue -ne "") {
            $eCell.Value = [double]$yearValue
        }
    }
}
